$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "301.32"
Set-TextValue $ws.Range("E2") "-0.73%"
Set-TextValue $ws.Range("G2") "22"

Set-TextValue $ws.Range("D3") "31.49"
Set-TextValue $ws.Range("E3") "-1.80%"
Set-TextValue $ws.Range("G3") "22"

Set-TextValue $ws.Range("D4") "5.090"
Set-TextValue $ws.Range("E4") "-3.12%"
Set-TextValue $ws.Range("G4") "22"

Set-TextValue $ws.Range("D5") "0.07363"
Set-TextValue $ws.Range("E5") "-2.25%"
Set-TextValue $ws.Range("G5") "22"

Set-TextValue $ws.Range("D6") "2.214"
Set-TextValue $ws.Range("E6") "44.98%"
Set-TextValue $ws.Range("G6") "22"

Set-TextValue $ws.Range("D7") "7.931"
Set-TextValue $ws.Range("E7") "0.23%"
Set-TextValue $ws.Range("G7") "22"

Set-TextValue $ws.Range("D8") "3.785"
Set-TextValue $ws.Range("E8") "-0.75%"
Set-TextValue $ws.Range("G8") "22"

Set-TextValue $ws.Range("D9") "0.9197"
Set-TextValue $ws.Range("E9") "-0.33%"
Set-TextValue $ws.Range("G9") "22"

Set-TextValue $ws.Range("D10") "0.1706"
Set-TextValue $ws.Range("G10") "22"

Set-TextValue $ws.Range("D11") "0.07598"
Set-TextValue $ws.Range("E11") "-3.61%"
Set-TextValue $ws.Range("G11") "22"

Set-TextValue $ws.Range("D12") "0.08140"
Set-TextValue $ws.Range("E12") "1.44%"
Set-TextValue $ws.Range("G12") "22"

Set-TextValue $ws.Range("E13") "-0.70%"
Set-TextValue $ws.Range("G13") "22"

Set-TextValue $ws.Range("D14") "0.09939"
Set-TextValue $ws.Range("E14") "0.19%"
Set-TextValue $ws.Range("G14") "22"

Set-TextValue $ws.Range("D15") "0.001493"
Set-TextValue $ws.Range("E15") "0.02%"
Set-TextValue $ws.Range("G15") "22"

Set-TextValue $ws.Range("D16") "0.006098"
Set-TextValue $ws.Range("E16") "-3.84%"
Set-TextValue $ws.Range("G16") "22"

Set-TextValue $ws.Range("D17") "3.463"
Set-TextValue $ws.Range("E17") "0.02%"
Set-TextValue $ws.Range("G17") "22"

Set-TextValue $ws.Range("D18") "2.225"
Set-TextValue $ws.Range("E18") "-0.33%"
Set-TextValue $ws.Range("G18") "22"

Set-TextValue $ws.Range("E19") "0.20%"
Set-TextValue $ws.Range("G19") "22"

Set-TextValue $ws.Range("D20") "0.1337"
Set-TextValue $ws.Range("E20") "1.61%"
Set-TextValue $ws.Range("G20") "22"

Set-TextValue $ws.Range("D21") "4.653"
Set-TextValue $ws.Range("E21") "3.87%"
Set-TextValue $ws.Range("G21") "22"

Set-TextValue $ws.Range("D22") "0.04644"
Set-TextValue $ws.Range("E22") "0.91%"
Set-TextValue $ws.Range("G22") "22"

Set-TextValue $ws.Range("E23") "-3.18%"
Set-TextValue $ws.Range("G23") "22"

Set-TextValue $ws.Range("E24") "0.96%"
Set-TextValue $ws.Range("G24") "22"

Set-TextValue $ws.Range("D25") "0.004481"
Set-TextValue $ws.Range("E25") "0.53%"
Set-TextValue $ws.Range("G25") "22"

Set-TextValue $ws.Range("E26") "-7.16%"
Set-TextValue $ws.Range("G26") "22"

Set-TextValue $ws.Range("E27") "49.97%"
Set-TextValue $ws.Range("G27") "22"

Set-TextValue $ws.Range("G28") "22"

Set-TextValue $ws.Range("G29") "22"

Set-TextValue $ws.Range("G30") "22"

Set-TextValue $ws.Range("G31") "22"

Set-TextValue $ws.Range("G32") "22"

Set-TextValue $ws.Range("G33") "22"

Set-TextValue $ws.Range("G34") "22"

Set-TextValue $ws.Range("G35") "22"

Set-TextValue $ws.Range("G36") "22"

Set-TextValue $ws.Range("G37") "22"

Set-TextValue $ws.Range("G38") "22"

Set-TextValue $ws.Range("D39") "0.01734"
Set-TextValue $ws.Range("E39") "2.09%"
Set-TextValue $ws.Range("G39") "22"

Set-TextValue $ws.Range("D40") "0.04529"
Set-TextValue $ws.Range("E40") "0.96%"
Set-TextValue $ws.Range("G40") "22"

Set-TextValue $ws.Range("D41") "0.007229"
Set-TextValue $ws.Range("E41") "4.17%"
Set-TextValue $ws.Range("G41") "22"

Set-TextValue $ws.Range("D42") "0.1347"
Set-TextValue $ws.Range("E42") "-0.37%"
Set-TextValue $ws.Range("G42") "22"

Set-TextValue $ws.Range("E43") "6.68%"
Set-TextValue $ws.Range("G43") "22"

Set-TextValue $ws.Range("E44") "-16.80%"
Set-TextValue $ws.Range("G44") "22"

Set-TextValue $ws.Range("D45") "0.00006298"
Set-TextValue $ws.Range("E45") "2.41%"
Set-TextValue $ws.Range("G45") "22"

Set-TextValue $ws.Range("E46") "-28.53%"
Set-TextValue $ws.Range("G46") "22"

Set-TextValue $ws.Range("D47") "0.8212"
Set-TextValue $ws.Range("E47") "-55.92%"
Set-TextValue $ws.Range("G47") "22"

Set-TextValue $ws.Range("G48") "22"

Set-TextValue $ws.Range("G49") "22"

Set-TextValue $ws.Range("G50") "22"

Set-TextValue $ws.Range("G51") "22"
